$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8483999999999999
$ws.Range("C2").Value = 0.6602
$ws.Range("B3").Value = 0.786
$ws.Range("C3").Value = 0.5570000000000001
$ws.Range("B4").Value = 0.2061366274507724
$ws.Range("C4").Value = 0.1589860871374152
